$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21664.25
$ws.Range("J3").Value = 21664.25
$ws.Range("L3").Value = 21664.25
$ws.Range("N3").Value = -21892.25
$ws.Range("H29").Value = 2057.1428
$ws.Range("J29").Value = 2760
$ws.Range("L29").Value = 8280
$ws.Range("N29").Value = -8842
$ws.Range("H38").Value = 785.7
$ws.Range("I38").Value = 317.44446
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 952.33338
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -580.33338
$ws.Range("N38").Value = -15744
$ws.Range("H40").Value = 1941.5
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 2049.6667
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 2049.6667
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = -2399.6667
$ws.Range("H58").Value = 675.44446
$ws.Range("H98").Value = 6203.5
$ws.Range("I98").Value = 6203.5
$ws.Range("K98").Value = 6203.5
$ws.Range("M98").Value = -4705.5
$ws.Range("H102").Value = 21664.25
$ws.Range("J102").Value = 21664.25
$ws.Range("L102").Value = 21664.25
$ws.Range("N102").Value = -28154.25
$ws.Range("H113").Value = 2832.2222
$ws.Range("I113").Value = 2570
$ws.Range("K113").Value = 2570
$ws.Range("M113").Value = 684
$ws.Range("H122").Value = 6203.5
$ws.Range("I122").Value = 6203.5
$ws.Range("K122").Value = 18610.5
$ws.Range("M122").Value = -16160.5
$ws.Range("H125").Value = 3179
$ws.Range("I125").Value = 2797.5
$ws.Range("K125").Value = 25177.5
$ws.Range("M125").Value = -22717.5
$ws.Range("H132").Value = 5468239.5
$ws.Range("I132").Value = 7577188
$ws.Range("J132").Value = 9783.529
$ws.Range("K132").Value = 22731564
$ws.Range("L132").Value = 29350.587
$ws.Range("M132").Value = -22729034
$ws.Range("N132").Value = -34410.587
$ws.Range("H135").Value = 1214.4286
$ws.Range("I135").Value = 506.6
$ws.Range("J135").Value = 2984
$ws.Range("K135").Value = 4559.400000000001
$ws.Range("L135").Value = 26856
$ws.Range("M135").Value = -2024.400000000001
$ws.Range("N135").Value = -31926
$ws.Range("H137").Value = 1032.7711
$ws.Range("I137").Value = 859.3488
$ws.Range("J137").Value = 1219.2
$ws.Range("K137").Value = 2578.0464
$ws.Range("L137").Value = 3657.6
$ws.Range("M137").Value = -28.04640000000018
$ws.Range("N137").Value = -8757.6
$ws.Range("H138").Value = 1381.82
$ws.Range("I138").Value = 629.04
$ws.Range("J138").Value = 1632.7467
$ws.Range("K138").Value = 1887.12
$ws.Range("L138").Value = 4898.2401
$ws.Range("M138").Value = 3252.88
$ws.Range("N138").Value = -15178.2401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1102.14
$ws.Range("I61").Value = 970.2195
$ws.Range("K61").Value = 970.2195
$ws.Range("M61").Value = -758.2195
$ws.Range("H102").Value = 55556856
$ws.Range("I102").Value = 55556856
$ws.Range("K102").Value = 55556856
$ws.Range("M102").Value = -55555234
$ws.Range("H122").Value = 2476.875
$ws.Range("I122").Value = 2476.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7430.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4980.625
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 3496.682
$ws.Range("I132").Value = 5069.1113
$ws.Range("J132").Value = 2408.077
$ws.Range("K132").Value = 15207.3339
$ws.Range("L132").Value = 7224.231000000001
$ws.Range("M132").Value = -12677.3339
$ws.Range("N132").Value = -12284.231
$ws.Range("H136").Value = 1102.14
$ws.Range("I136").Value = 970.2195
$ws.Range("K136").Value = 2910.6585
$ws.Range("M136").Value = -360.6585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 200
$ws.Range("J8").Value = 200
$ws.Range("L8").Value = 200
$ws.Range("N8").Value = -480
$ws.Range("H94").Value = 20834776
$ws.Range("I94").Value = 25001310
$ws.Range("K94").Value = 25001310
$ws.Range("M94").Value = -25000859
$ws.Range("H134").Value = 4436.814
$ws.Range("I134").Value = 1500.6562
$ws.Range("K134").Value = 4501.9686
$ws.Range("M134").Value = -1966.9686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66668148
$ws.Range("I16").Value = 76924480
$ws.Range("K16").Value = 76924480
$ws.Range("M16").Value = -76924193
$ws.Range("H86").Value = 2159811.2
$ws.Range("I86").Value = 3177118.5
$ws.Range("K86").Value = 3177118.5
$ws.Range("M86").Value = -3175995.5
$ws.Range("H89").Value = 2159811.2
$ws.Range("I89").Value = 3177118.5
$ws.Range("K89").Value = 15885592.5
$ws.Range("M89").Value = -15879976.5
$ws.Range("H94").Value = 606.375
$ws.Range("I94").Value = 528
$ws.Range("J94").Value = 632.5
$ws.Range("K94").Value = 528
$ws.Range("L94").Value = 632.5
$ws.Range("M94").Value = -77
$ws.Range("N94").Value = -1534.5
$ws.Range("H99").Value = 2438.5334
$ws.Range("I99").Value = 2314.8333
$ws.Range("J99").Value = 2933.3333
$ws.Range("K99").Value = 2314.8333
$ws.Range("L99").Value = 2933.3333
$ws.Range("M99").Value = -816.8332999999998
$ws.Range("N99").Value = -5929.3333
$ws.Range("H107").Value = 452.1579
$ws.Range("I107").Value = 394.69232
$ws.Range("J107").Value = 576.6667
$ws.Range("K107").Value = 394.69232
$ws.Range("L107").Value = 576.6667
$ws.Range("M107").Value = 1525.30768
$ws.Range("N107").Value = -4416.6667
$ws.Range("H113").Value = 66668148
$ws.Range("I113").Value = 76924480
$ws.Range("K113").Value = 76924480
$ws.Range("M113").Value = -76922310
$ws.Range("H122").Value = 1460.5834
$ws.Range("I122").Value = 1043.8182
$ws.Range("J122").Value = 1813.2307
$ws.Range("K122").Value = 3131.4546
$ws.Range("L122").Value = 5439.6921
$ws.Range("M122").Value = -681.4546
$ws.Range("N122").Value = -10339.6921
$ws.Range("H126").Value = 2438.5334
$ws.Range("I126").Value = 2314.8333
$ws.Range("J126").Value = 2933.3333
$ws.Range("K126").Value = 6944.499899999999
$ws.Range("L126").Value = 8799.999899999999
$ws.Range("M126").Value = -4474.499899999999
$ws.Range("N126").Value = -13739.9999
$ws.Range("H132").Value = 6079.68
$ws.Range("I132").Value = 6333.048
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 18999.144
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -16469.144
$ws.Range("N132").Value = -19308.5
$ws.Range("H134").Value = 856
$ws.Range("I134").Value = 833.3273
$ws.Range("K134").Value = 2499.9819
$ws.Range("M134").Value = 35.01809999999978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5655.2104
$ws.Range("J107").Value = 7533.0713
$ws.Range("L107").Value = 22599.2139
$ws.Range("N107").Value = -26439.2139
$ws.Range("H131").Value = 32259562
$ws.Range("I131").Value = 111111440
$ws.Range("J131").Value = 1976.2273
$ws.Range("K131").Value = 333334320
$ws.Range("L131").Value = 5928.6819
$ws.Range("M131").Value = -333329280
$ws.Range("N131").Value = -16008.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 10000
$ws.Range("N48").Value = -10970
$ws.Range("H70").Value = 4602.25
$ws.Range("I70").Value = 4700
$ws.Range("J70").Value = 4569.6665
$ws.Range("K70").Value = 4700
$ws.Range("L70").Value = 4569.6665
$ws.Range("M70").Value = -4430
$ws.Range("N70").Value = -5109.6665
$ws.Range("H73").Value = 4602.25
$ws.Range("I73").Value = 4700
$ws.Range("J73").Value = 4569.6665
$ws.Range("K73").Value = 4700
$ws.Range("L73").Value = 4569.6665
$ws.Range("M73").Value = -3764
$ws.Range("N73").Value = -6441.6665
$ws.Range("H126").Value = 1832.6
$ws.Range("I126").Value = 1635.3334
$ws.Range("J126").Value = 2128.5
$ws.Range("K126").Value = 4906.0002
$ws.Range("L126").Value = 6385.5
$ws.Range("M126").Value = -2436.0002
$ws.Range("N126").Value = -11325.5
$ws.Range("H132").Value = 2390.1765
$ws.Range("I132").Value = 2071.2307
$ws.Range("J132").Value = 3426.75
$ws.Range("K132").Value = 6213.6921
$ws.Range("L132").Value = 10280.25
$ws.Range("M132").Value = -3683.6921
$ws.Range("N132").Value = -15340.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 335.25
$ws.Range("I55").Value = 262.94446
$ws.Range("J55").Value = 428.2143
$ws.Range("K55").Value = 262.94446
$ws.Range("L55").Value = 428.2143
$ws.Range("M55").Value = -89.94445999999999
$ws.Range("N55").Value = -774.2143
$ws.Range("H61").Value = 1404.2
$ws.Range("I61").Value = 1130.25
$ws.Range("K61").Value = 1130.25
$ws.Range("M61").Value = -928.25
$ws.Range("H113").Value = 1404.2
$ws.Range("I113").Value = 1130.25
$ws.Range("K113").Value = 1130.25
$ws.Range("M113").Value = 1039.75
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H132").Value = 29308.406
$ws.Range("I132").Value = 1960.2
$ws.Range("J132").Value = 47954.91
$ws.Range("K132").Value = 5880.6
$ws.Range("L132").Value = 143864.73
$ws.Range("M132").Value = -3350.6
$ws.Range("N132").Value = -148924.73
$ws.Range("H136").Value = 3704.027
$ws.Range("I136").Value = 4236.4136
$ws.Range("K136").Value = 12709.2408
$ws.Range("M136").Value = -10159.2408

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H96").Value = 3268
$ws.Range("I96").Value = 3399.9
$ws.Range("K96").Value = 3399.9
$ws.Range("M96").Value = -2026.9
$ws.Range("H122").Value = 138890600
$ws.Range("I122").Value = 138890600
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 416671800
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -416669350
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 3182.8
$ws.Range("I132").Value = 3308.68
$ws.Range("J132").Value = 2553.4
$ws.Range("K132").Value = 9926.039999999999
$ws.Range("L132").Value = 7660.200000000001
$ws.Range("M132").Value = -7396.039999999999
$ws.Range("N132").Value = -12720.2
